$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C->D, old D->E)
$ws.Columns.Item(3).Insert()

# New column header
$ws.Range("C1").Value = "Book Value [ExclRevalReserve]/Share (Rs.)"

# New column C values (Book Value [ExclRevalReserve]/Share (Rs.)) per row.
# These look numeric but must be stored as text, matching the rest of the
# sheet (which uses inline/shared text strings, not numbers). Temporarily
# format the range as Text, assign the values, then restore the style so no
# residual number formatting is left behind on the cells.
$values = @{
    2  = "213.71"
    3  = "207.49"
    4  = "201.85"
    5  = "37.18"
    6  = "35.46"
    7  = "32.75"
    8  = "30.05"
    9  = "29.07"
    10 = "17.21"
    11 = "15.15"
    12 = "12.36"
    13 = "16.25"
    14 = "12.31"
    15 = "11.84"
    16 = "9.45"
    17 = "6.61"
    18 = "12.34"
    19 = "10.47"
    20 = "9.50"
}

$dataRange = $ws.Range("C2:C20")
$dataRange.NumberFormat = "@"

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

$dataRange.Style = "Normal"
